$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.732.70"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "2.526.57"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.32"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.85"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "2.526.66"
$ws.Range("E10").Value = "  +6.14%  "
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "2.988.76"
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("E15").Value = "  +3.39%  "
$ws.Range("D16").Value = "68.594.83"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "2.509.77"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.51"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.26"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.83"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.25"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("E26").Value = "  -5.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.03"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "2.690.82"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "518.72"
$ws.Range("E30").Value = "  +4.13%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0894"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.82"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.05"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.46"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.69"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  +5.54%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.70"
$ws.Range("E46").Value = "  +7.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.56"
$ws.Range("E47").Value = "  +2.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.522"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("E51").Value = "  -0.09%  "

# Remove the temporary text-number formatting so cells retain default (no explicit) style
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()

